# This script "re-derives" the N (seas_id_y -> season_ending_year_y) and
# O (season_ending_year_y -> player_id_y) columns of the mip_winners sheet,
# and relabels the related header cells (K, L, N, O) to match the renamed
# columns coming out of the upstream merge (…_x / …_y suffix rework).
#
# Before:
#   K = seas_id_x        L = player_id         N = seas_id_y              O = season_ending_year_y
# After:
#   K = seas_id          L = player_id_x       N = season_ending_year_y   O = player_id_y
#
# The N column used to hold a raw seas_id number and the O column held the
# season-ending-year text; after the merge rework N now holds the
# season-ending-year (as text, taken from the old O column) while O holds a
# brand new player_id number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header relabel (row 1) ---------------------------------------------
$ws.Range("K1").Value = "seas_id"
$ws.Range("L1").Value = "player_id_x"
$ws.Range("N1").Value = "season_ending_year_y"
$ws.Range("O1").Value = "player_id_y"

# --- New player_id_y values (column O) for rows 2..40 -------------------
$playerIdY = @{
    2  = 5025
    3  = 3251
    4  = 2223
    5  = 2936
    6  = 609
    7  = 3966
    8  = 5061
    9  = 1947
    10 = 718
    11 = 2638
    12 = 1970
    13 = 3994
    14 = 4452
    15 = 3104
    16 = 13
    17 = 1142
    18 = 2091
    19 = 3772
    20 = 582
    21 = 566
    22 = 5265
    23 = 1951
    24 = 2516
    25 = 4941
    26 = 2321
    27 = 1173
    28 = 87
    29 = 2174
    30 = 1945
    31 = 1122
    32 = 1461
    33 = 3404
    34 = 4038
    35 = 4525
    36 = 4411
    37 = 3097
    38 = 3088
    39 = 1085
    40 = 148
}

# Column N must keep storing TEXT (it used to hold the season-ending-year
# string that lived in column O) rather than being re-inferred as a number,
# so force a text number format on the range before writing the values, then
# drop back to the Normal style once the text values are committed.
$colN = $ws.Range("N2:N40")
$colN.NumberFormat = "@"

for ($row = 2; $row -le 40; $row++) {
    # season_ending_year_y (column N) is the calendar year the season ended,
    # decreasing by one every row starting at 2024 on row 2.
    $seasonEndingYear = 2026 - $row

    # Column N: was a bare seas_id number, now the season-ending-year text
    # (same value that used to live in column O).
    $ws.Cells.Item($row, 14).Value = [string]$seasonEndingYear

    # Column O: was the season-ending-year text, now a brand new
    # player_id_y number.
    $ws.Cells.Item($row, 15).Value = $playerIdY[$row]
}

$colN.Style = "Normal"
